$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-10 19:20:33", 13, 12, 4, 3, 5, 0, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.9230769230769231),
    @("2023-12-10 21:05:03", 11, 6, 1, 0, 5, 0, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.5454545454545454),
    @("2023-12-10 22:30:25", 25, 19, 6, 2, 11, 0, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.76)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
